$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Append the 9 new "GVAO_PLCC*" indicator rows (295-303).
#    Values are written column-by-column (all of column B, then all
#    of column C, then D formulas, then E, then F) so that the new
#    shared-string entries are appended to sharedStrings.xml in the
#    same order as in the target workbook:
#      GVAO_PLCCACT, GVAO_PLCCNSW, GVAO_PLCCQSL, GVAO_PLCCSAL,
#      GVAO_PLCCTAS, GVAO_PLCCVIC, GVAO_PLCCNTY, GVAO_PLCCWAL,
#      GVAO_PLCCAUS, GVAO_PLCC
# -----------------------------------------------------------------

$mnemonics = @("GVAO_PLCCACT","GVAO_PLCCNSW","GVAO_PLCCQSL","GVAO_PLCCSAL","GVAO_PLCCTAS","GVAO_PLCCVIC","GVAO_PLCCNTY","GVAO_PLCCWAL","GVAO_PLCCAUS")
$firstRow = 295

for ($i = 0; $i -lt $mnemonics.Length; $i++) {
    $r = $firstRow + $i
    $ws.Range("A$r").Value = "AID"
    $ws.Range("B$r").Value = $mnemonics[$i]
    $ws.Range("C$r").Value = "STATES"
}

for ($i = 0; $i -lt $mnemonics.Length; $i++) {
    $r = $firstRow + $i
    $ws.Range("D$r").Formula = '=B' + $r + '&"_"&C' + $r
}

for ($i = 0; $i -lt $mnemonics.Length; $i++) {
    $r = $firstRow + $i
    $ws.Range("E$r").Value = "GVAO_PLCC"
    $ws.Range("F$r").Value = "Sum"
}

# -----------------------------------------------------------------
# 2. Re-apply the AutoFilter over the new A1:F294 range, filtering
#    column B (Mnemonic) to the GVAOLCC* state series and column C
#    (Division) to STATES. This recreates the <autoFilter> element
#    and automatically hides every row that does not satisfy both
#    criteria (matching the target workbook's hidden rows).
# -----------------------------------------------------------------

$ws.AutoFilterMode = $false

$ws.Range("A1:F294").AutoFilter(2, @("GVAOLCCACT","GVAOLCCAUS","GVAOLCCNSW","GVAOLCCNTY","GVAOLCCQSL","GVAOLCCSAL","GVAOLCCTAS","GVAOLCCVIC","GVAOLCCWAL"), 7)
$ws.Range("A1:F294").AutoFilter(3, @("STATES"), 7)

# -----------------------------------------------------------------
# 3. Update the workbook-level _FilterDatabase defined name so it
#    refers to the new filtered range.
# -----------------------------------------------------------------

for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name -eq "Sheet1!_FilterDatabase") {
        $nm.RefersTo = "=Sheet1!`$A`$1:`$F`$294"
    }
}

# -----------------------------------------------------------------
# 4. Update the sheet view: scroll back to the top-left and select
#    cell B235 (matching the saved selection in the target file).
# -----------------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B235").Select()
